$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append new row 13 with the 2025-12-07 data.
# Column A holds a date-like string that must stay plain text (matching the
# existing rows), so force a text number format, assign the value, then clear
# the formatting so no style index is left behind on the cell.
$cellA = $ws.Cells.Item(13, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "12/07/2025"
$cellA.ClearFormats()

$ws.Cells.Item(13, 2).Value = 13193.07
$ws.Cells.Item(13, 3).Value = 0.1724867103395756
$ws.Cells.Item(13, 4).Value = 0.8275132896604244
$ws.Cells.Item(13, 5).Value = -94.31
$ws.Cells.Item(13, 6).Value = -21.68
$ws.Cells.Item(13, 7).Value = -19574.32
$ws.Cells.Item(13, 8).Value = -64.2
$ws.Cells.Item(13, 9).Value = -527.22
$ws.Cells.Item(13, 10).Value = -18.81
